# "adicionando versao final da interface e corrigindo enderecos de
# carregamento dos modelos"
#
# Rows 5-10 get overwritten with new respondents / scores, and rows
# 11-23 are brand new respondents appended below the original data.
# Most respondents only answered question 1 (column B), so columns
# C:F are blanked out (kept as empty text, matching a quote-prefixed
# blank entry) while G (Media Geral) mirrors the single answer. Two
# respondents ("gabriel" / row 14 and "t" / row 17) answered all 5
# questions, so their C:F + G stay fully populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# name, Q1, Q2, Q3, Q4, Q5, Media Geral  ($null = no answer for that question)
$rows = @(
    @(5,  "Murilo Silva",   3.19, $null, $null, $null, $null, 3.19),
    @(6,  "Murilo Silva",   2.9,  $null, $null, $null, $null, 2.9),
    @(7,  "Murilo Silva",   2.14, $null, $null, $null, $null, 2.14),
    @(8,  "teste",          0.52, $null, $null, $null, $null, 0.52),
    @(9,  "a",               0.59, $null, $null, $null, $null, 0.59),
    @(10, "b",               1.99, $null, $null, $null, $null, 1.99),
    @(11, "Charles",         1.36, $null, $null, $null, $null, 1.36),
    @(12, "Murilo",          3.45, $null, $null, $null, $null, 3.45),
    @(13, "Murilo alves",    2.79, $null, $null, $null, $null, 2.79),
    @(14, "gabriel",         2.7,  2.71,  2.62,  2.75,  2.6,   2.68),
    @(15, "testando",        0.52, $null, $null, $null, $null, 0.52),
    @(16, "Destro",          4.49, $null, $null, $null, $null, 4.49),
    @(17, "t",                1.91, 1.9,   3.28,  1.14,  2.38,  2.12),
    @(18, "a",                1.36, $null, $null, $null, $null, 1.36),
    @(19, "a",                2.25, $null, $null, $null, $null, 2.25),
    @(20, "a",                2.11, $null, $null, $null, $null, 2.11),
    @(21, "a",                1.36, $null, $null, $null, $null, 1.36),
    @(22, "a",                0.99, $null, $null, $null, $null, 0.99),
    @(23, "testando0000",     1.21, $null, $null, $null, $null, 1.21)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    for ($col = 2; $col -le 6; $col++) {
        $answer = $entry[$col]
        if ($answer -eq $null) {
            # No answer for this question -> blank text cell (keeps the
            # column typed as text instead of a fully-empty/number cell).
            $ws.Cells.Item($r, $col).Value = "'"
            $ws.Cells.Item($r, $col).Style = "Normal"
        } else {
            $ws.Cells.Item($r, $col).Value = $answer
        }
    }
    $ws.Cells.Item($r, 7).Value = $entry[7]
}
